$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "RAWData"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "Bytes"

$ws.Activate()
$excel.ActiveWindow.Zoom = 145
$ws.Range("T305").Select()

Write-Host "done"
